# chore: update Sheets via scheduled runner
# Refresh cached market-board pricing / leve-profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) for the affected leve
# rows across each job sheet, as produced by the scheduled pricing-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 500000060
$ws.Range("I9").Value = 500000060
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 500000060
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -499999891
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("ALC")  # row 40
$ws.Range("H40").Value = 1866.1915
$ws.Range("I40").Value = 1867.6316
$ws.Range("K40").Value = 1867.6316
$ws.Range("M40").Value = -1692.6316

$ws = $wb.Worksheets.Item("ALC")  # row 41
$ws.Range("H41").Value = 784.125
$ws.Range("I41").Value = 671.875
$ws.Range("J41").Value = 896.375
$ws.Range("K41").Value = 671.875
$ws.Range("L41").Value = 896.375
$ws.Range("M41").Value = -231.875
$ws.Range("N41").Value = -1776.375

$ws = $wb.Worksheets.Item("ALC")  # row 51
$ws.Range("H51").Value = 2625

$ws = $wb.Worksheets.Item("ALC")  # row 55
$ws.Range("H55").Value = 647.1818
$ws.Range("I55").Value = 818.625
$ws.Range("J55").Value = 190
$ws.Range("K55").Value = 818.625
$ws.Range("L55").Value = 190
$ws.Range("M55").Value = -604.625
$ws.Range("N55").Value = -618

$ws = $wb.Worksheets.Item("ALC")  # row 92
$ws.Range("H92").Value = 12963886
$ws.Range("I92").Value = 2646153.5
$ws.Range("J92").Value = 37038596
$ws.Range("K92").Value = 2646153.5
$ws.Range("L92").Value = 37038596
$ws.Range("M92").Value = -2644905.5
$ws.Range("N92").Value = -37041092

$ws = $wb.Worksheets.Item("ALC")  # row 101
$ws.Range("H101").Value = 1759.0667
$ws.Range("J101").Value = 3491.4285
$ws.Range("L101").Value = 10474.2855
$ws.Range("N101").Value = -13718.2855

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 1659.7709
$ws.Range("I137").Value = 1132.7727
$ws.Range("J137").Value = 2105.6924
$ws.Range("K137").Value = 3398.3181
$ws.Range("L137").Value = 6317.0772
$ws.Range("M137").Value = -848.3181
$ws.Range("N137").Value = -11417.0772

$ws = $wb.Worksheets.Item("ALC")  # row 138
$ws.Range("H138").Value = 3350.2615
$ws.Range("I138").Value = 1565.9131
$ws.Range("J138").Value = 4327.405
$ws.Range("K138").Value = 4697.7393
$ws.Range("L138").Value = 12982.215
$ws.Range("M138").Value = 442.2606999999998
$ws.Range("N138").Value = -23262.215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5028.8374
$ws.Range("I32").Value = 4667.154
$ws.Range("J32").Value = 6596.1333
$ws.Range("K32").Value = 4667.154
$ws.Range("L32").Value = 6596.1333
$ws.Range("M32").Value = -4380.154
$ws.Range("N32").Value = -7170.1333

$ws = $wb.Worksheets.Item("ARM")  # row 45
$ws.Range("H45").Value = 4922.654
$ws.Range("I45").Value = 9048.416999999999
$ws.Range("J45").Value = 1386.2858
$ws.Range("K45").Value = 9048.416999999999
$ws.Range("L45").Value = 1386.2858
$ws.Range("M45").Value = -8671.416999999999
$ws.Range("N45").Value = -2140.2858

$ws = $wb.Worksheets.Item("ARM")  # row 61
$ws.Range("H61").Value = 3533.72
$ws.Range("I61").Value = 4541.2573
$ws.Range("J61").Value = 1182.8
$ws.Range("K61").Value = 4541.2573
$ws.Range("L61").Value = 1182.8
$ws.Range("M61").Value = -4329.2573
$ws.Range("N61").Value = -1606.8

$ws = $wb.Worksheets.Item("ARM")  # row 102
$ws.Range("H102").Value = 4632707
$ws.Range("I102").Value = 9261486
$ws.Range("K102").Value = 9261486
$ws.Range("M102").Value = -9259864

$ws = $wb.Worksheets.Item("ARM")  # row 136
$ws.Range("H136").Value = 3533.72
$ws.Range("I136").Value = 4541.2573
$ws.Range("J136").Value = 1182.8
$ws.Range("K136").Value = 13623.7719
$ws.Range("L136").Value = 3548.4
$ws.Range("M136").Value = -11073.7719
$ws.Range("N136").Value = -8648.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9436.387000000001
$ws.Range("I20").Value = 1227.9
$ws.Range("J20").Value = 24360.908
$ws.Range("K20").Value = 1227.9
$ws.Range("L20").Value = 24360.908
$ws.Range("M20").Value = -980.9000000000001
$ws.Range("N20").Value = -24854.908

$ws = $wb.Worksheets.Item("BSM")  # row 107
$ws.Range("H107").Value = 1422.2222
$ws.Range("I107").Value = 1414.2858
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 1414.2858
$ws.Range("L107").Value = 1450
$ws.Range("M107").Value = 505.7141999999999
$ws.Range("N107").Value = -5290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2264.7932
$ws.Range("I31").Value = 1497.7142
$ws.Range("J31").Value = 2508.8635
$ws.Range("K31").Value = 1497.7142
$ws.Range("L31").Value = 2508.8635
$ws.Range("M31").Value = -1202.7142
$ws.Range("N31").Value = -3098.8635

$ws = $wb.Worksheets.Item("CRP")  # row 34
$ws.Range("H34").Value = 2264.7932
$ws.Range("I34").Value = 1497.7142
$ws.Range("J34").Value = 2508.8635
$ws.Range("K34").Value = 1497.7142
$ws.Range("L34").Value = 2508.8635
$ws.Range("M34").Value = -1295.7142
$ws.Range("N34").Value = -2912.8635

$ws = $wb.Worksheets.Item("CRP")  # row 82
$ws.Range("H82").Value = 32000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws = $wb.Worksheets.Item("CRP")  # row 85
$ws.Range("H85").Value = 32000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")  # row 134
$ws.Range("H134").Value = 3145.16
$ws.Range("I134").Value = 3322.3914
$ws.Range("J134").Value = 1107
$ws.Range("K134").Value = 9967.174199999999
$ws.Range("L134").Value = 3321
$ws.Range("M134").Value = -7432.174199999999
$ws.Range("N134").Value = -8391

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 42571.43
$ws.Range("J37").Value = 42571.43
$ws.Range("L37").Value = 127714.29
$ws.Range("N37").Value = -127938.29

$ws = $wb.Worksheets.Item("CUL")  # row 68
$ws.Range("H68").Value = 2859.7917
$ws.Range("I68").Value = 4024.5312
$ws.Range("J68").Value = 1928
$ws.Range("K68").Value = 12073.5936
$ws.Range("L68").Value = 5784
$ws.Range("M68").Value = -11262.5936
$ws.Range("N68").Value = -7406

$ws = $wb.Worksheets.Item("CUL")  # row 71
$ws.Range("H71").Value = 2859.7917
$ws.Range("I71").Value = 4024.5312
$ws.Range("J71").Value = 1928
$ws.Range("K71").Value = 36220.7808
$ws.Range("L71").Value = 17352
$ws.Range("M71").Value = -32164.7808
$ws.Range("N71").Value = -25464

$ws = $wb.Worksheets.Item("CUL")  # row 107
$ws.Range("H107").Value = 1164.434
$ws.Range("J107").Value = 1308.9111
$ws.Range("L107").Value = 3926.7333
$ws.Range("N107").Value = -7766.7333

$ws = $wb.Worksheets.Item("CUL")  # row 109
$ws.Range("H109").Value = 2194.5
$ws.Range("I109").Value = 837.1667
$ws.Range("J109").Value = 3212.5
$ws.Range("K109").Value = 2511.5001
$ws.Range("L109").Value = 9637.5
$ws.Range("M109").Value = -1471.5001
$ws.Range("N109").Value = -11717.5

$ws = $wb.Worksheets.Item("CUL")  # row 113
$ws.Range("H113").Value = 1071915.9
$ws.Range("I113").Value = 1429018.8
$ws.Range("J113").Value = 476744.38
$ws.Range("K113").Value = 4287056.4
$ws.Range("L113").Value = 1430233.14
$ws.Range("M113").Value = -4284886.4
$ws.Range("N113").Value = -1434573.14

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 318.63635
$ws.Range("I2").Value = 336.2
$ws.Range("J2").Value = 143
$ws.Range("K2").Value = 336.2
$ws.Range("L2").Value = 143
$ws.Range("M2").Value = -223.2
$ws.Range("N2").Value = -369

$ws = $wb.Worksheets.Item("GSM")  # row 113
$ws.Range("H113").Value = 66667932
$ws.Range("I113").Value = 111112090
$ws.Range("J113").Value = 1700.5
$ws.Range("K113").Value = 111112090
$ws.Range("L113").Value = 1700.5
$ws.Range("M113").Value = -111109920
$ws.Range("N113").Value = -6040.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 6666.3335
$ws.Range("J12").Value = 6666.3335
$ws.Range("L12").Value = 6666.3335
$ws.Range("N12").Value = -7006.3335

$ws = $wb.Worksheets.Item("LTW")  # row 21
$ws.Range("H21").Value = 52503.5
$ws.Range("J21").Value = 52503.5
$ws.Range("L21").Value = 52503.5
$ws.Range("N21").Value = -52851.5

$ws = $wb.Worksheets.Item("LTW")  # row 22
$ws.Range("H22").Value = 18519910
$ws.Range("I22").Value = 111111110
$ws.Range("J22").Value = 1670
$ws.Range("K22").Value = 111111110
$ws.Range("L22").Value = 1670
$ws.Range("M22").Value = -111110815
$ws.Range("N22").Value = -2260

$ws = $wb.Worksheets.Item("LTW")  # row 27
$ws.Range("H27").Value = 18519910
$ws.Range("I27").Value = 111111110
$ws.Range("J27").Value = 1670
$ws.Range("K27").Value = 111111110
$ws.Range("L27").Value = 1670
$ws.Range("M27").Value = -111111003
$ws.Range("N27").Value = -1884

$ws = $wb.Worksheets.Item("LTW")  # row 46
$ws.Range("H46").Value = 22223362
$ws.Range("I46").Value = 37037980
$ws.Range("J46").Value = 1433.3334
$ws.Range("K46").Value = 37037980
$ws.Range("L46").Value = 1433.3334
$ws.Range("M46").Value = -37037792
$ws.Range("N46").Value = -1809.3334

$ws = $wb.Worksheets.Item("LTW")  # row 55
$ws.Range("H55").Value = 16129359
$ws.Range("I55").Value = 272.2143
$ws.Range("J55").Value = 29412136
$ws.Range("K55").Value = 272.2143
$ws.Range("L55").Value = 29412136
$ws.Range("M55").Value = -99.21429999999998
$ws.Range("N55").Value = -29412482

$ws = $wb.Worksheets.Item("LTW")  # row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7142.3335
$ws.Range("J15").Value = 7142.3335
$ws.Range("L15").Value = 7142.3335
$ws.Range("N15").Value = -7718.3335
